$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1. First paragraph: append a red "(This is a change ... )" annotation,
#    split across three runs (mirrors how Word would leave separate
#    edit-session runs with identical red formatting).
# -----------------------------------------------------------------------
$d.Content.Find.Execute(
    "This is a Microsoft word document.", $true, $false, $false, $false,
    $false, $true, 1, $false, "This is a Microsoft word document.  ", 2
) | Out-Null

$p1 = $d.Paragraphs(1)
$insPoint = $p1.Range.End - 1

$chunk1 = "(This is a change " + [string]([char]0x2013) + " Ve"
$r1 = $d.Range($insPoint, $insPoint)
$r1.InsertAfter($chunk1)
$r1.Font.Color = 192

$insPoint2 = $insPoint + $chunk1.Length
$chunk2 = "rsion for branch alternate"
$r2 = $d.Range($insPoint2, $insPoint2)
$r2.InsertAfter($chunk2)
$r2.Font.Color = 192

$insPoint3 = $insPoint2 + $chunk2.Length
$chunk3 = ")"
$r3 = $d.Range($insPoint3, $insPoint3)
$r3.InsertAfter($chunk3)
$r3.Font.Color = 192

# -----------------------------------------------------------------------
# 2. "Crispian's Day speech from Shakespear's Henry V [Source - Wikipedia]"
#    paragraph: no text changes, just re-flow run/proofErr boundaries.
#    Doing in-place Find/Replace over the touched spans causes identically
#    formatted neighbouring runs to coalesce, matching the target shape.
# -----------------------------------------------------------------------
$d.Content.Find.Execute(
    " Day speech from ", $true, $false, $false, $false, $false, $true, 1,
    $false, " Day speech from ", 2
) | Out-Null

$d.Content.Find.Execute(
    "Henry V [Source", $true, $false, $false, $false, $false, $true, 1,
    $false, "Henry V [Source", 2
) | Out-Null

# -----------------------------------------------------------------------
# 3. Append two new paragraphs at the end of the document: one using the
#    "larger" style with explicit shading/spacing, and a final bare one.
#    Building the raw paragraph XML (rather than InsertParagraphAfter)
#    avoids inheriting the previous paragraph-mark run formatting.
# -----------------------------------------------------------------------
$endOfDoc = $d.Content.End
$tailRange = $d.Range($endOfDoc, $endOfDoc)
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$newParasXml =
    '<w:p xmlns:w="' + $wNs + '">' +
        '<w:pPr>' +
            '<w:pStyle w:val="larger"/>' +
            '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
            '<w:spacing w:before="0" w:beforeAutospacing="0" w:after="150" w:afterAutospacing="0"/>' +
        '</w:pPr>' +
    '</w:p>' +
    '<w:p xmlns:w="' + $wNs + '"/>'
$tailRange.InsertXML($newParasXml)

# -----------------------------------------------------------------------
# 4. Styles cleanup: drop the now-unused "apple-converted-space" and
#    "Hyperlink" character styles (delete highest index first so the
#    style collection stays valid).
# -----------------------------------------------------------------------
$d.Styles("Hyperlink").Delete()
$d.Styles("apple-converted-space").Delete()

Write-Output "edit complete"
